$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 ("metricas_recorrencia_anual") with refreshed Bibi data
$ws.Range("C8").Value = 934
$ws.Range("D8").Value = 155
$ws.Range("E8").Value = 779
$ws.Range("F8").Value = 6.357670221493026
$ws.Range("G8").Value = 83.40471092077088
$ws.Range("H8").Value = 16.59528907922912
